# Work on the input format for the Outstation Patcher.
# Adds a "Comment" column (AL) value to several rows of the Work_List sheet,
# including a brand-new free-text comment on row 4, and updates the
# sheet's selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work_List")

# Populate the new "Comment" column (AL) entries that are missing.
$ws.Range("AL2").Value  = "Comment"
$ws.Range("AL4").Value  = "Zzzzzzzz"
$ws.Range("AL6").Value  = "Comment"
$ws.Range("AL10").Value = "Comment"
$ws.Range("AL11").Value = "Comment"

# Move the selection to reflect where the user ended up after editing.
$ws.Range("AL12").Select()
